# Update Strategy Script and .gitignore File.
#
# 1) Change the custom date / date-time number formats from using
#    "yyyy-mm-dd" style separators to "yyyy/mm/dd" style separators.
# 2) Move the sheet's active cell / selection from A9 to H25.
# 3) Update the electricity_rate values in G24:G26 from 54.2 to 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Number format changes -------------------------------------------
# Column A (date) uses "yyyy-mm-dd"; columns D/E (datetime_start/end) use
# "yyyy-mm-dd hh:mm:ss". Re-apply the same cells with slash separators.
$ws.Range("A3:A26").NumberFormat = "yyyy/mm/dd"
$ws.Range("D3:E26").NumberFormat = "yyyy/mm/dd hh:mm:ss"

# --- 2) Update electricity_rate values for the last three rows ----------
$ws.Range("G24:G26").Value = 51

# --- 3) Move the active selection to H25 ---------------------------------
$ws.Range("H25").Select() | Out-Null
